$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column I: duplicate column H values for rows 1-4
$ws.Range('I1').Value = 'Girl,People,School,Animals,Cat,Fish,Lion,Elephant,Sheep,Chair,Door,Jacket,City,Beautiful,Yellow,Orange,Rainbow,Ice'
$ws.Range('I2').Value = 'she,he,them,i sit on the bench,she sits on the bench,moon,money,noodles,orange,oven,Pineapple,park ,queen,question,rainbow,time,umbrella,ear,zebra,violin,white,hear,glue,paint,five apples,eight apples'
$ws.Range('I3').Value = 'ribhi'
$ws.Range('I4').Value = 'ribhi'

# Column J: move the old dictation marker (previously in column I) to column J
$ws.Range('J1').Value = 'dictation'
$ws.Range('J2').Value = 'dictation'
$ws.Range('J3').Value = 'dictation'
$ws.Range('J4').Value = 'dictation'

# New row 5 (a new multiple-choice question entry)
$ws.Range('A5').Value = 'A1'
$ws.Range('B5').Value = 'Listening'
$ws.Range('C5').Value = 'Unit1'
$ws.Range('D5').Value = 'In this lesson you will hear some words , try typing them correctly . You can hear each word how much ever you like . All of these words are taken before so you wont find them so hard (if you have taken notes ;)​'
$ws.Range('E5').Value = 'ستسمع في هذا الدرس بعض الكلمات، حاول كتابتها بشكل صحيح. يمكنك سماع كل كلمة كم مرة تريد. تم أخذ كل هذه الكلمات لذا لن تجدها صعبة للغاية (إذا كنت قد قمت بتدوين الملاحظات:)​'
$ws.Range('F5').Value = 'Write down what you are hearing'
$ws.Range('G5').Value = 'اكتب ما تسمعه'
$ws.Range('H5').Value = 'Girl,People,School'
$ws.Range('I5').Value = 'Girl '
$ws.Range('J5').Value = 'multipleChoice'

# Column widths for new columns I and J
$ws.Columns.Item(9).ColumnWidth = 13.498697916666666
$ws.Columns.Item(10).ColumnWidth = 17.330729166666668

# Update selection to match target view state
$ws.Range('J16').Select()
